$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.490.25"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.98%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.324.45"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.85%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.03"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.80%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.586"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.320.73"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.74%  "
$ws.Range("E10").Value = "  -4.09%  "
$ws.Range("E11").Value = "  -2.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "45.37"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.68%  "
$ws.Range("E13").Value = "  -4.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "663.93"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.82%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.866.89"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.64%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.38"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.87%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.682.57"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.81%  "
$ws.Range("E18").Value = "  -1.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.320.24"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.55%  "
$ws.Range("E21").Value = "  -2.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.887"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.42"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.99"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.52%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "98.63"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.41%  "
$ws.Range("E26").Value = "  -6.46%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.67"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.63%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "33.72"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.23"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.38"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.51%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.63%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "591.72"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.53%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "10.91"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.04%  "
$ws.Range("E34").Value = "  -1.93%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.698.30"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -8.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "56.79"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.24"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -14.38%  "
$ws.Range("E39").Value = "  +0.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "33.42"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.62"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.85%  "
$ws.Range("E42").Value = "  -7.04%  "
$ws.Range("E43").Value = "  -3.62%  "
$ws.Range("E44").Value = "  -6.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.23"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0405"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.85%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.59"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.11%  "
$ws.Range("E48").Value = "  -2.10%  "
$ws.Range("E50").Value = "  -3.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "126.72"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.11%  "
